$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.485.25"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "1.919.20"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4828"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4081"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08227"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.023"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").Value = "1.922.44"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.051"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06806"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001041"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "29.515.70"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.651"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.196"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").Value = "2.166.38"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.639"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.117"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.024"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09578"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.543"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.563"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.384"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02288"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06139"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.183"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5984"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.042"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1855"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.415"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("E45").Value = "  +2.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07604"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5586"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.89%  "
$ws.Range("E50").Value = "  +4.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.42%  "
